$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Item #2 - 0603 resistor CRG0603F300R): more resistors added to designator list
$ws.Range("B3").Value = "R10,R12,R8,R13,R15,R16,R17"
$ws.Range("C3").Value = 7

# Row 7 (Item #6 - LED APD3224SURCK-F01): more LEDs added to designator list
$ws.Range("B7").Value = "D4,D2,D1,D3,D5,D6,D7"
$ws.Range("C7").Value = 7

# Row 12 (Item #11 - IC2): part swapped from Onsemi MC74ACT05DG to TI SN74HCS125QBQARQ1
$ws.Range("D12").Value = "Texas Instruments"
$ws.Range("E12").Value = "SN74HCS125QBQARQ1"
$ws.Range("F12").Value = "Automotive Schmitt-trigger inputs quadruple bus buffer gates with 3-state outputs 14-WQFN -40 to 125 "

# Column width adjustments (B widened, F widened)
$ws.Columns(2).ColumnWidth = 25
$ws.Columns(6).ColumnWidth = 84.6667

# Update selection / scroll position
$ws.Range("A4").Select()
